# Updates the cached "datetimeFigureOut" field text (date placeholder) on the
# slide master and every slide layout from 05-06-2014 to 18-06-2014, and
# touches the "Scared" mood-rectangle text run on slide 1 (re-affirms its
# text / formatting), matching the upstream template2.pptx commit.

$p = $ppt.ActivePresentation

function Update-DateText($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.HasTextFrame) {
            if ($sh.TextFrame.HasText) {
                if ($sh.TextFrame.TextRange.Text -eq "05-06-2014") {
                    $sh.TextFrame.TextRange.Text = "18-06-2014"
                }
            }
        }
    }
}

# Slide master's own date placeholder.
$master = $p.SlideMaster
Update-DateText $master.Shapes

# Every slide layout's date placeholder.
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DateText $layout.Shapes
}

# Slide 1: re-affirm the "Scared" rectangle's text run.
$slide = $p.Slides.Item(1)
$scared = $slide.Shapes.Item("Rectangle 12")
$scared.TextFrame.TextRange.Text = "Scared"
